# Update leve-profit calculation cells (columns H-N) across several sheets
# to reflect refreshed market-board pricing data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 932371.4  # H6 was 1004053.75
$ws.Cells.Item(6, 10).Value = 2973  # J6 was 3797.3333
$ws.Cells.Item(6, 12).Value = 8919  # L6 was 11391.9999
$ws.Cells.Item(6, 14).Value = -9143  # N6 was -11615.9999
$ws.Cells.Item(51, 8).Value = 19268  # H51 was 17614.428
$ws.Cells.Item(51, 9).Value = 52000  # I51 was 27100
$ws.Cells.Item(51, 10).Value = 2902  # J51 was 4967
$ws.Cells.Item(51, 11).Value = 52000  # K51 was 27100
$ws.Cells.Item(51, 12).Value = 2902  # L51 was 4967
$ws.Cells.Item(51, 13).Value = -51516  # M51 was -26616
$ws.Cells.Item(51, 14).Value = -3870  # N51 was -5935
$ws.Cells.Item(103, 8).Value = 549.6177  # H103 was 554.08826
$ws.Cells.Item(103, 9).Value = 401.16666  # I103 was 411.8
$ws.Cells.Item(103, 10).Value = 581.4286  # J103 was 578.62067
$ws.Cells.Item(103, 11).Value = 1203.49998  # K103 was 1235.4
$ws.Cells.Item(103, 12).Value = 1744.2858  # L103 was 1735.86201
$ws.Cells.Item(103, 13).Value = -617.4999800000001  # M103 was -649.4000000000001
$ws.Cells.Item(103, 14).Value = -2916.2858  # N103 was -2907.86201
$ws.Cells.Item(111, 8).Value = 10012802  # H111 was 7702622
$ws.Cells.Item(111, 9).Value = 52029  # I111 was 22023.4
$ws.Cells.Item(111, 11).Value = 156087  # K111 was 66070.20000000001
$ws.Cells.Item(111, 13).Value = -153020  # M111 was -63003.20000000001
$ws.Cells.Item(129, 8).Value = 884.1177  # H129 was 888.5238000000001
$ws.Cells.Item(129, 10).Value = 955.6  # J129 was 1015.7059
$ws.Cells.Item(129, 12).Value = 2866.8  # L129 was 3047.1177
$ws.Cells.Item(129, 14).Value = -12866.8  # N129 was -13047.1177
$ws.Cells.Item(137, 8).Value = 2506.2856  # H137 was 2749.25
$ws.Cells.Item(137, 9).Value = 1747.0834  # I137 was 1798.5
$ws.Cells.Item(137, 10).Value = 3518.5557  # J137 was 4333.8335
$ws.Cells.Item(137, 11).Value = 5241.2502  # K137 was 5395.5
$ws.Cells.Item(137, 12).Value = 10555.6671  # L137 was 13001.5005
$ws.Cells.Item(137, 13).Value = -2691.2502  # M137 was -2845.5
$ws.Cells.Item(137, 14).Value = -15655.6671  # N137 was -18101.5005
$ws.Cells.Item(138, 8).Value = 3789.7646  # H138 was 4396.271
$ws.Cells.Item(138, 9).Value = 1296.871  # I138 was 1445.2963
$ws.Cells.Item(138, 10).Value = 5878.4053  # J138 was 8190.381
$ws.Cells.Item(138, 11).Value = 3890.613  # K138 was 4335.8889
$ws.Cells.Item(138, 12).Value = 17635.2159  # L138 was 24571.143
$ws.Cells.Item(138, 13).Value = 1249.387  # M138 was 804.1111000000001
$ws.Cells.Item(138, 14).Value = -27915.2159  # N138 was -34851.143
$ws.Cells.Item(141, 8).Value = 2755.4167  # H141 was 3260.4614
$ws.Cells.Item(141, 9).Value = 2755.4167  # I141 was 2962.3635
$ws.Cells.Item(141, 10).Value = 0  # J141 was 4900
$ws.Cells.Item(141, 11).Value = 8266.250100000001  # K141 was 8887.0905
$ws.Cells.Item(141, 12).Value = 0  # L141 was 14700
$ws.Cells.Item(141, 13).Value = -3086.250100000001  # M141 was -3707.0905
$ws.Cells.Item(141, 14).Value = $null  # N141 was -25060
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(8, 8).Value = 0  # H8 was 3000
$ws.Cells.Item(8, 9).Value = 0  # I8 was 3000
$ws.Cells.Item(8, 11).Value = 0  # K8 was 3000
$ws.Cells.Item(8, 13).Value = $null  # M8 was -2856
$ws.Cells.Item(45, 8).Value = 3828.2104  # H45 was 3907.5
$ws.Cells.Item(45, 9).Value = 3288.9092  # I45 was 3297.818
$ws.Cells.Item(45, 10).Value = 4569.75  # J45 was 4865.5713
$ws.Cells.Item(45, 11).Value = 3288.9092  # K45 was 3297.818
$ws.Cells.Item(45, 12).Value = 4569.75  # L45 was 4865.5713
$ws.Cells.Item(45, 13).Value = -2911.9092  # M45 was -2920.818
$ws.Cells.Item(45, 14).Value = -5323.75  # N45 was -5619.5713
$ws.Cells.Item(106, 8).Value = 41890  # H106 was 43990
$ws.Cells.Item(106, 10).Value = 41890  # J106 was 43990
$ws.Cells.Item(106, 12).Value = 41890  # L106 was 43990
$ws.Cells.Item(106, 14).Value = -44414  # N106 was -46514
$ws.Cells.Item(107, 8).Value = 0  # H107 was 27400.5
$ws.Cells.Item(107, 10).Value = 0  # J107 was 27400.5
$ws.Cells.Item(107, 12).Value = 0  # L107 was 27400.5
$ws.Cells.Item(107, 14).Value = $null  # N107 was -35080.5
$ws.Cells.Item(108, 8).Value = 40571  # H108 was 42400
$ws.Cells.Item(108, 10).Value = 40571  # J108 was 42400
$ws.Cells.Item(108, 12).Value = 40571  # L108 was 42400
$ws.Cells.Item(108, 14).Value = -48251  # N108 was -50080
$ws.Cells.Item(109, 8).Value = 34800  # H109 was 0
$ws.Cells.Item(109, 10).Value = 34800  # J109 was 0
$ws.Cells.Item(109, 12).Value = 34800  # L109 was 0
$ws.Cells.Item(109, 14).Value = -37574  # N109 was empty (new cell)
$ws.Cells.Item(112, 8).Value = 0  # H112 was 15236.2
$ws.Cells.Item(112, 10).Value = 0  # J112 was 15236.2
$ws.Cells.Item(112, 12).Value = 0  # L112 was 15236.2
$ws.Cells.Item(112, 14).Value = $null  # N112 was -18190.2
$ws.Cells.Item(118, 8).Value = 38728.09  # H118 was 38980
$ws.Cells.Item(118, 10).Value = 38728.09  # J118 was 38980
$ws.Cells.Item(118, 12).Value = 38728.09  # L118 was 38980
$ws.Cells.Item(118, 14).Value = -42042.09  # N118 was -42294
$ws.Cells.Item(119, 8).Value = 40849  # H119 was 37844.668
$ws.Cells.Item(119, 10).Value = 40849  # J119 was 37844.668
$ws.Cells.Item(119, 12).Value = 40849  # L119 was 37844.668
$ws.Cells.Item(119, 14).Value = -50525  # N119 was -47520.668
$ws.Cells.Item(132, 8).Value = 2976.84  # H132 was 3265.422
$ws.Cells.Item(132, 9).Value = 2871.4546  # I132 was 3190.923
$ws.Cells.Item(132, 11).Value = 8614.363799999999  # K132 was 9572.769
$ws.Cells.Item(132, 13).Value = -6084.363799999999  # M132 was -7042.769
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7, 8).Value = 0  # H7 was 1326.6666
$ws.Cells.Item(7, 10).Value = 0  # J7 was 1326.6666
$ws.Cells.Item(7, 12).Value = 0  # L7 was 1326.6666
$ws.Cells.Item(7, 14).Value = $null  # N7 was -1552.6666
$ws.Cells.Item(99, 8).Value = 1299.9474  # H99 was 1215.2632
$ws.Cells.Item(99, 9).Value = 951.4286  # I99 was 900
$ws.Cells.Item(99, 10).Value = 2275.8  # J99 was 2397.5
$ws.Cells.Item(99, 11).Value = 951.4286  # K99 was 900
$ws.Cells.Item(99, 12).Value = 2275.8  # L99 was 2397.5
$ws.Cells.Item(99, 13).Value = 546.5714  # M99 was 598
$ws.Cells.Item(99, 14).Value = -5271.8  # N99 was -5393.5
$ws.Cells.Item(134, 8).Value = 2137.4  # H134 was 2261.6382
$ws.Cells.Item(134, 9).Value = 2057.122  # I134 was 2194.2896
$ws.Cells.Item(134, 10).Value = 2503.111  # J134 was 2546
$ws.Cells.Item(134, 11).Value = 6171.366  # K134 was 6582.8688
$ws.Cells.Item(134, 12).Value = 7509.333  # L134 was 7638
$ws.Cells.Item(134, 13).Value = -3636.366  # M134 was -4047.8688
$ws.Cells.Item(134, 14).Value = -12579.333  # N134 was -12708
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1342.8572  # H16 was 1225
$ws.Cells.Item(16, 9).Value = 1200  # I16 was 887.5
$ws.Cells.Item(16, 11).Value = 1200  # K16 was 887.5
$ws.Cells.Item(16, 13).Value = -913  # M16 was -600.5
$ws.Cells.Item(31, 8).Value = 32803.453  # H31 was 35562.133
$ws.Cells.Item(31, 9).Value = 44703.74  # I31 was 39647.54
$ws.Cells.Item(31, 10).Value = 5432.8  # J31 was 9007
$ws.Cells.Item(31, 11).Value = 44703.74  # K31 was 39647.54
$ws.Cells.Item(31, 12).Value = 5432.8  # L31 was 9007
$ws.Cells.Item(31, 13).Value = -44408.74  # M31 was -39352.54
$ws.Cells.Item(31, 14).Value = -6022.8  # N31 was -9597
$ws.Cells.Item(34, 8).Value = 32803.453  # H34 was 35562.133
$ws.Cells.Item(34, 9).Value = 44703.74  # I34 was 39647.54
$ws.Cells.Item(34, 10).Value = 5432.8  # J34 was 9007
$ws.Cells.Item(34, 11).Value = 44703.74  # K34 was 39647.54
$ws.Cells.Item(34, 12).Value = 5432.8  # L34 was 9007
$ws.Cells.Item(34, 13).Value = -44501.74  # M34 was -39445.54
$ws.Cells.Item(34, 14).Value = -5836.8  # N34 was -9411
$ws.Cells.Item(58, 8).Value = 11369.92  # H58 was 10527.63
$ws.Cells.Item(58, 9).Value = 1037.909  # I58 was 969
$ws.Cells.Item(58, 10).Value = 87138  # J58 was 65489.75
$ws.Cells.Item(58, 11).Value = 1037.909  # K58 was 969
$ws.Cells.Item(58, 12).Value = 87138  # L58 was 65489.75
$ws.Cells.Item(58, 13).Value = -834.9090000000001  # M58 was -766
$ws.Cells.Item(58, 14).Value = -87544  # N58 was -65895.75
$ws.Cells.Item(74, 8).Value = 22410.908  # H74 was 22649.8
$ws.Cells.Item(74, 10).Value = 22410.908  # J74 was 22649.8
$ws.Cells.Item(74, 12).Value = 22410.908  # L74 was 22649.8
$ws.Cells.Item(74, 14).Value = -24158.908  # N74 was -24397.8
$ws.Cells.Item(77, 8).Value = 22410.908  # H77 was 22649.8
$ws.Cells.Item(77, 10).Value = 22410.908  # J77 was 22649.8
$ws.Cells.Item(77, 12).Value = 67232.724  # L77 was 67949.39999999999
$ws.Cells.Item(77, 14).Value = -75968.724  # N77 was -76685.39999999999
$ws.Cells.Item(113, 8).Value = 1342.8572  # H113 was 1225
$ws.Cells.Item(113, 9).Value = 1200  # I113 was 887.5
$ws.Cells.Item(113, 11).Value = 1200  # K113 was 887.5
$ws.Cells.Item(113, 13).Value = 970  # M113 was 1282.5
$ws.Cells.Item(136, 8).Value = 11369.92  # H136 was 10527.63
$ws.Cells.Item(136, 9).Value = 1037.909  # I136 was 969
$ws.Cells.Item(136, 10).Value = 87138  # J136 was 65489.75
$ws.Cells.Item(136, 11).Value = 3113.727  # K136 was 2907
$ws.Cells.Item(136, 12).Value = 261414  # L136 was 196469.25
$ws.Cells.Item(136, 13).Value = -563.7270000000003  # M136 was -357
$ws.Cells.Item(136, 14).Value = -266514  # N136 was -201569.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 828.23  # H131 was 824.23
$ws.Cells.Item(131, 10).Value = 840.25  # J131 was 836.0833
$ws.Cells.Item(131, 12).Value = 2520.75  # L131 was 2508.2499
$ws.Cells.Item(131, 14).Value = -12600.75  # N131 was -12588.2499
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 1264200  # H5 was 1684433.4
$ws.Cells.Item(5, 10).Value = 18933.334  # J5 was 26650
$ws.Cells.Item(5, 12).Value = 18933.334  # L5 was 26650
$ws.Cells.Item(5, 14).Value = -19157.334  # N5 was -26874
$ws.Cells.Item(109, 8).Value = 0  # H109 was 7356.6665
$ws.Cells.Item(109, 10).Value = 0  # J109 was 7356.6665
$ws.Cells.Item(109, 12).Value = 0  # L109 was 7356.6665
$ws.Cells.Item(109, 14).Value = $null  # N109 was -9436.666499999999
$ws.Cells.Item(113, 8).Value = 1624  # H113 was 1373.1428
$ws.Cells.Item(113, 10).Value = 1684.7142  # J113 was 1402.1666
$ws.Cells.Item(113, 12).Value = 1684.7142  # L113 was 1402.1666
$ws.Cells.Item(113, 14).Value = -6024.7142  # N113 was -5742.1666
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 3486.7097  # H132 was 3586.0334
$ws.Cells.Item(132, 9).Value = 4090.2727  # I132 was 4089.9546
$ws.Cells.Item(132, 10).Value = 2011.3334  # J132 was 2200.25
$ws.Cells.Item(132, 11).Value = 12270.8181  # K132 was 12269.8638
$ws.Cells.Item(132, 12).Value = 6034.0002  # L132 was 6600.75
$ws.Cells.Item(132, 13).Value = -9740.8181  # M132 was -9739.863799999999
$ws.Cells.Item(132, 14).Value = -11094.0002  # N132 was -11660.75
$ws.Cells.Item(136, 8).Value = 2178.0833  # H136 was 2087.8
$ws.Cells.Item(136, 9).Value = 2104.2222  # I136 was 1995.3182
$ws.Cells.Item(136, 10).Value = 2399.6667  # J136 was 2766
$ws.Cells.Item(136, 11).Value = 6312.6666  # K136 was 5985.9546
$ws.Cells.Item(136, 12).Value = 7199.000100000001  # L136 was 8298
$ws.Cells.Item(136, 13).Value = -3762.6666  # M136 was -3435.9546
$ws.Cells.Item(136, 14).Value = -12299.0001  # N136 was -13398
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 77092.664  # H46 was 100429
$ws.Cells.Item(46, 10).Value = 77092.664  # J46 was 100429
$ws.Cells.Item(46, 12).Value = 77092.664  # L46 was 100429
$ws.Cells.Item(46, 14).Value = -77554.664  # N46 was -100891
$ws.Cells.Item(132, 8).Value = 2143.0293  # H132 was 2324.9333
$ws.Cells.Item(132, 9).Value = 1928  # I132 was 2035.762
$ws.Cells.Item(132, 10).Value = 2537.25  # J132 was 2999.6667
$ws.Cells.Item(132, 11).Value = 5784  # K132 was 6107.286
$ws.Cells.Item(132, 12).Value = 7611.75  # L132 was 8999.000100000001
$ws.Cells.Item(132, 13).Value = -3254  # M132 was -3577.286
$ws.Cells.Item(132, 14).Value = -12671.75  # N132 was -14059.0001
$ws.Cells.Item(134, 8).Value = 77092.664  # H134 was 100429
$ws.Cells.Item(134, 10).Value = 77092.664  # J134 was 100429
$ws.Cells.Item(134, 12).Value = 231277.992  # L134 was 301287
$ws.Cells.Item(134, 14).Value = -236347.992  # N134 was -306357
